$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a literal text value into a cell without letting Excel's
# automatic type-inference turn numeric-looking strings (e.g. "1.00",
# "353.85") into real numbers, and without leaving any NumberFormat /
# quote-prefix style residue behind on the cell (matches the original
# inlineStr text cells, which all carry the default/no style).
function Set-TextCell($sheet, $ref, $val) {
    $c = $sheet.Range($ref)
    $c.Value = "'" + $val
    $c.ClearFormats()
}

Set-TextCell $ws "D2" "52.134.85"
Set-TextCell $ws "E2" "  +0.71%  "

Set-TextCell $ws "D3" "2.903.14"
Set-TextCell $ws "E3" "  +3.58%  "

Set-TextCell $ws "D5" "353.85"
Set-TextCell $ws "E5" "  +0.20%  "

Set-TextCell $ws "D6" "113.44"
Set-TextCell $ws "E6" "  +1.13%  "

Set-TextCell $ws "D7" "0.556"
Set-TextCell $ws "E7" "  -0.42%  "

Set-TextCell $ws "E8" "  +0.04%  "

Set-TextCell $ws "D9" "0.623"
Set-TextCell $ws "E9" "  +0.04%  "

Set-TextCell $ws "D10" "39.59"
Set-TextCell $ws "E10" "  -1.43%  "

Set-TextCell $ws "D11" "0.0874"
Set-TextCell $ws "E11" "  +4.31%  "

Set-TextCell $ws "E12" "  +0.78%  "

Set-TextCell $ws "D13" "19.77"
Set-TextCell $ws "E13" "  -0.63%  "

Set-TextCell $ws "D14" "7.73"
Set-TextCell $ws "E14" "  -0.48%  "

Set-TextCell $ws "D15" "3.366.23"
Set-TextCell $ws "E15" "  +3.78%  "

Set-TextCell $ws "D16" "2.909.11"
Set-TextCell $ws "E16" "  +3.20%  "

Set-TextCell $ws "D17" "0.984"
Set-TextCell $ws "E17" "  +3.06%  "

Set-TextCell $ws "D18" "52.221.22"
Set-TextCell $ws "E18" "  +0.83%  "

Set-TextCell $ws "E19" "  +1.54%  "

Set-TextCell $ws "D20" "7.59"
Set-TextCell $ws "E20" "  -0.17%  "

Set-TextCell $ws "D21" "14.07"
Set-TextCell $ws "E21" "  +3.77%  "

Set-TextCell $ws "D22" "0.0₃0979"
Set-TextCell $ws "E22" "  +0.63%  "

Set-TextCell $ws "D23" "71.01"
Set-TextCell $ws "E23" "  +1.02%  "

Set-TextCell $ws "D24" "269.49"
Set-TextCell $ws "E24" "  +0.74%  "

Set-TextCell $ws "D25" "2.81"
Set-TextCell $ws "E25" "  +1.83%  "

Set-TextCell $ws "D26" "0.181"
Set-TextCell $ws "E26" "  +12.23%  "

Set-TextCell $ws "D27" "26.74"
Set-TextCell $ws "E27" "  +2.32%  "

Set-TextCell $ws "D28" "1.00"
Set-TextCell $ws "E28" "  -0.03%  "

Set-TextCell $ws "D29" "10.64"
Set-TextCell $ws "E29" "  +2.29%  "

Set-TextCell $ws "E30" "  +15.02%  "

Set-TextCell $ws "D31" "6.80"
Set-TextCell $ws "E31" "  +11.33%  "

Set-TextCell $ws "B32" "Toncoin"
Set-TextCell $ws "C32" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws "D32" "2.27"
Set-TextCell $ws "E32" "  -0.66%  "

Set-TextCell $ws "B33" "InjectiveProtocol"
Set-TextCell $ws "C33" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws "D33" "37.41"
Set-TextCell $ws "E33" "  -4.45%  "

Set-TextCell $ws "D34" "6.09"
Set-TextCell $ws "E34" "  +10.31%  "

Set-TextCell $ws "D35" "53.03"
Set-TextCell $ws "E35" "  +1.46%  "

Set-TextCell $ws "D36" "0.0451"
Set-TextCell $ws "E36" "  +0.36%  "

Set-TextCell $ws "E37" "  -0.08%  "

Set-TextCell $ws "E38" "  +4.70%  "

Set-TextCell $ws "D39" "18.79"
Set-TextCell $ws "E39" "  -1.06%  "

Set-TextCell $ws "E40" "  +1.45%  "

Set-TextCell $ws "D41" "2.72"
Set-TextCell $ws "E41" "  +8.08%  "

Set-TextCell $ws "E42" "  +1.09%  "

Set-TextCell $ws "D43" "23.08"
Set-TextCell $ws "E43" "  +5.19%  "

Set-TextCell $ws "B44" "WEMIXToken"
Set-TextCell $ws "C44" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws "D44" "2.18"
Set-TextCell $ws "E44" "  -2.32%  "

Set-TextCell $ws "B45" "Monero"
Set-TextCell $ws "C45" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws "D45" "117.86"
Set-TextCell $ws "E45" "  -1.60%  "

Set-TextCell $ws "E46" "  +1.86%  "

Set-TextCell $ws "D47" "3.52"
Set-TextCell $ws "E47" "  -0.31%  "

Set-TextCell $ws "D48" "2.180.55"
Set-TextCell $ws "E48" "  +2.93%  "

Set-TextCell $ws "D49" "0.259"
Set-TextCell $ws "E49" "  +17.73%  "

Set-TextCell $ws "D50" "0.0355"
Set-TextCell $ws "E50" "  +12.41%  "

Set-TextCell $ws "D51" "0.951"
Set-TextCell $ws "E51" "  -2.99%  "

Write-Host "done"
